$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the existing backup codes with newly generated ones.
$ws.Range("A2").Value = "SZ5TCJ42NKEA"
$ws.Range("A3").Value = "BYH6BYARV86T"
$ws.Range("A4").Value = "YK1J90E88BG6"

# Append a new block of backup codes (Crdc TC02) below the existing data.
$ws.Range("A13").Value = "CWNRQHM7CMQ3"
$ws.Range("A14").Value = "PQYKR4S0SCAE"
$ws.Range("A15").Value = "HB2TNX3R72X2"
$ws.Range("A16").Value = "NP11MCW57D3V"

# Move the active selection to match the saved workbook state.
$ws.Range("C6").Select()
